$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 19.42991633333333
$ws.Range("N2").Value = 58.289749
$ws.Range("O2").Value = 0.08673502554925175
$ws.Range("P2").Value = 0.08673502554925173
$ws.Range("Q2").Value = 0.1230820433328889
$ws.Range("R2").Value = 1.107738389996
$ws.Range("S2").Value = 0.08673502554925175
$ws.Range("T2").Value = 0.08673502554925173

# Row 3
$ws.Range("O3").Value = 0.04103322570207864
$ws.Range("P3").Value = 0.04103322570207864
$ws.Range("S3").Value = 0.04103322570207864
$ws.Range("T3").Value = 0.04103322570207864

# Row 4
$ws.Range("M4").Value = 105.042315
$ws.Range("N4").Value = 315.126945
$ws.Range("O4").Value = 0.4689082402093144
$ws.Range("P4").Value = 0.4689082402093144
$ws.Range("Q4").Value = 0.6654080514199999
$ws.Range("R4").Value = 5.988672462779999
$ws.Range("S4").Value = 0.4689082402093144
$ws.Range("T4").Value = 0.4689082402093144

# Row 5
$ws.Range("M5").Value = 2.834125333333334
$ws.Range("N5").Value = 8.502376000000002
$ws.Range("O5").Value = 0.01265151784388958
$ws.Range("P5").Value = 0.01265151784388958
$ws.Range("Q5").Value = 0.01795323927822223
$ws.Range("R5").Value = 0.161579153504
$ws.Range("S5").Value = 0.01265151784388958
$ws.Range("T5").Value = 0.01265151784388958

# Row 6
$ws.Range("M6").Value = 87.516249
$ws.Range("N6").Value = 262.548747
$ws.Range("O6").Value = 0.3906719906954657
$ws.Range("P6").Value = 0.3906719906954657
$ws.Range("Q6").Value = 0.554386265332
$ws.Range("R6").Value = 4.989476387988
$ws.Range("S6").Value = 0.3906719906954657
$ws.Range("T6").Value = 0.3906719906954657
